$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CBM thickness")
$ws.Activate()
$ws.Rows("9:9").Delete()
$ws.Range("A9").Select()
